# Applies the "Allen Time Interval" edit:
#  - Inserts two new rows (alphabetically sorted) for "funds" and
#    "funding_received" between "description" and "idText"
#  - Flips every value in column B from TRUE to FALSE
#  - Re-applies the sort over the (now larger) A1:A24 range so the
#    worksheet's sortState reflects the new extent
#  - Moves the active cell selection to F18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 12 ("idText"),
# pushing everything below down by two rows. Excel carries the B-column
# style/formatting down automatically.
$ws.Rows("12:13").Insert() | Out-Null

# Populate the two new label cells. Set A13 ("funds") before A12
# ("funding_received") so they are appended to the shared-strings table
# in the same order as the target workbook (funds, then
# funding_received).
$ws.Range("A13").Value() = "funds"
$ws.Range("A12").Value() = "funding_received"

# All boolean flags in column B become FALSE.
[void]($ws.Range("B1:B24").Value() = $false)

# Re-sort column A (header-less) over the new full range so the
# worksheet's stored sort state covers A1:A24 instead of the old A1:A21.
$sortRange = $ws.Range("A1:A24")
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("A1:A24")) | Out-Null
$ws.Sort.SetRange($sortRange) | Out-Null
$ws.Sort.Header = 0
$ws.Sort.Apply() | Out-Null

# Update the active selection to match the target workbook.
$ws.Range("F18").Select() | Out-Null
